$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.479.89"
$ws.Range("D3").Value = "3.672.65"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'648.87"
$ws.Range("E5").Value = "  -4.64%  "
$ws.Range("D6").Value = "'159.70"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'7.12"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D13").Value = "4.292.08"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "'32.56"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "3.662.36"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "69.463.84"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "'15.97"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "'464.84"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "'9.76"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "'0.644"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").Value = "'79.45"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "3.818.75"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D26").Value = "'0.0000125"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").Value = "'10.79"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'8.94"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").Value = "'2.61"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").Value = "'1.67"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'2.00"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.46"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.60"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").Value = "'0.164"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").Value = "3.662.51"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").Value = "'8.38"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D39").Value = "'5.90"
$ws.Range("E39").Value = "  -5.65%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'177.55"
$ws.Range("E41").Value = "  +4.27%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0895"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.18"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("D44").Value = "'0.926"
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").Value = "'46.64"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  -3.08%  "
$ws.Range("D48").Value = "'26.95"
$ws.Range("E48").Value = "  -5.08%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "'0.000266"
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.82"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  -5.70%  "
